$wb = $excel.ActiveWorkbook

# --- Sheet "Gaz" ---
$ws = $wb.Worksheets.Item("Gaz")

# Header row
$ws.Range("B1").Value = "Last Price"
$ws.Range("C1").Value = "Last Volume"
$ws.Range("D1").Value = "End of Day Index"

# Row 2 (2025-06-16) - new numeric values
$ws.Range("B2").Value = 37.15
$ws.Range("C2").Value = 13680
$ws.Range("D2").Value = 36.934

# Row 3 (2025-06-17) - shift old E3:G3 values into B3:D3
$ws.Range("B3").Value = 38.95
$ws.Range("C3").Value = 24000
$ws.Range("D3").Value = 38.201

# Clear old columns E, F, G (no longer used)
$ws.Range("E1:G3").Clear()

# --- Sheet "CO2" ---
$ws2 = $wb.Worksheets.Item("CO2")
$ws2.Range("B2").Value = 74.7
